$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the row for Conta 004550415 / Nome DIOGO / Saldo 3548.54 ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "004550415") {
        $ws.Rows.Item($r).Delete()
        break
    }
}

# --- Update Saldo for Conta 004272426 / Nome RODRIGO: -1450.77 -> -1441.68 ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "004272426") {
        $ws.Cells.Item($r, 3).Value = -1441.68
        break
    }
}
